$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 3 (shifts existing rows 3-10 down to 4-11),
# copying formatting from the row above so the new row matches the table style.
$ws.Rows.Item(3).Insert()

# Fill in the new recruit for S'23: Marketing Intern, Sumin Sung.
$ws.Range("B3").Value = "Marketing Intern"
$ws.Range("A3").Value = "Sumin Sung"

# Match formatting (style) of the other data rows (e.g. row 4) for the new row 3.
$ws.Range("A4:C4").Copy()
$ws.Range("A3:C3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update the selection to match where the edit happened.
$ws.Range("A3").Select()
